$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# -----------------------------------------------------------------------
$p1 = $d.Paragraphs(1)

$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Aurora Beast Hunter, a high variance 5-reel and 40-payline slot game with bonus rounds and free spins. Play Aurora Beast Hunter for free.</w:t></w:r></w:p>'

$metaText = "Meta description: Read our review of Aurora Beast Hunter, a high variance 5-reel and 40-payline slot game with bonus rounds and free spins. Play Aurora Beast Hunter for free."

$insertPoint = $d.Range($p1.Range.End, $p1.Range.End)
$insertPoint.InsertXML($metaXml)

# The inserted runs land inside the following paragraph (Heading2), so
# split right after the inserted text to give the new content its own
# paragraph.
$splitPos = $p1.Range.End + $metaText.Length
$splitRange = $d.Range($splitPos, $splitPos)
$splitRange.InsertParagraphAfter()

# The new paragraph inherited the Heading2 style from the split point;
# reset it back to the body-text (Normal) style used elsewhere.
$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

# -----------------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph that used to sit just
#    before the closing italic "meta description" paragraph.
# -----------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    $pText = $p.Range.Text.TrimEnd("`r")
    if ($pText -eq "Play Aurora Beast Hunter Free - Review of High Variance Slot Game" -and $i -ne 1) {
        $p.Range.Delete()
        break
    }
}

# -----------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the new
#    DALLE image-prompt copy, keeping its italic run formatting and
#    leading empty run intact.
# -----------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)

$newImageText = "Create a Cartoon Feature Image for Aurora Beast Hunter Slot Game DALLE, please create a colorful cartoon-style feature image for the slot game, Aurora Beast Hunter. The image should feature a happy Maya warrior with glasses. The Maya warrior should be standing in a desert canyon with an array of weapon symbols and laser guns around. In the background, there should be a beehive-shaped reel-set with some of the game's characters and symbols. The image should be bold, vibrant, and eye-catching to attract players to the game. Thank you!"

$pStart = $lastPara.Range.Start
$pEnd = $lastPara.Range.End
$textRange = $d.Range($pStart, $pEnd - 1)
$textRange.Text = $newImageText

Write-Output "done"
